# Corrected beamwidth for 18nm-S
# Rows 6-13 hold the 18nm-S measurements (A = "26T2300_18nm_S").
# Column H = sigma_B / um, Column I = sigma_B_err / % for that row's
# orientation (G = "vertical" -> 432 um / 8%, G = "horizontal" -> 382 um / 5%).
# All dependent formula cells (E, F, J, N, P, Q, R) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value  = 432
$ws.Range("I6").Value  = 0.08
$ws.Range("H7").Value  = 382
$ws.Range("I7").Value  = 0.05
$ws.Range("H8").Value  = 432
$ws.Range("I8").Value  = 0.08
$ws.Range("H9").Value  = 382
$ws.Range("I9").Value  = 0.05
$ws.Range("H10").Value = 432
$ws.Range("I10").Value = 0.08
$ws.Range("H11").Value = 382
$ws.Range("I11").Value = 0.05
$ws.Range("H12").Value = 432
$ws.Range("I12").Value = 0.08
$ws.Range("H13").Value = 382
$ws.Range("I13").Value = 0.05

# The workbook had a stale/unused external link to results.xlsx; remove it
# (this also drops the <externalReferences> entry from workbook.xml).
$wb.BreakLink("file:///C:\Users\Thomas\OneDrive\PhD\DP-Experiment\analysis\results.xlsx", 1) | Out-Null

# Move the active selection to I16, matching the saved cursor position.
$ws.Range("I16").Select() | Out-Null
